# Update input/output data on Hoja1:
#  - add the e-mail address as a new row (A2), backed by a new shared string
#  - size column A to fit the new content (bestFit width ~27 chars)
#  - leave the selection on E4, matching the refreshed workbook view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "contacto.diego.c@gmail.com"

# 26.14 characters of "ColumnWidth" rounds/stores as an internal column
# width of exactly 27 (Excel's char-width -> internal-unit conversion),
# matching the auto-fit width for this header/value pair.
$ws.Columns.Item(1).ColumnWidth = 26.14

$ws.Range("E4").Select()
